# Update cryptocurrency price/volume figures (inline-string cells) for the
# GitHub Actions scheduled refresh. D = Price, E = Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.165.86"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -1.96%  "
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.84"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.66%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.13"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("E6").Value = "  -0.69%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4223"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -2.14%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3679"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -1.75%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07229"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -1.61%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8539"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -2.96%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.94"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.96%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.33"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("E13").Value = "  -0.61%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07078"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -0.80%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.298"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -2.81%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.61"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("E17").Value = "  -0.88%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008839"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  -3.05%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.251.19"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -1.70%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.113"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("E23").Value = "  -2.65%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.052.40"
$ws.Range("D24").Style = $style
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.979"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("E26").Value = "  -2.09%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.195"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +2.85%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.36"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -1.08%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.231"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -2.81%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.24"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("E31").Value = "  -1.17%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.189"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -3.20%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7478"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -4.02%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.938"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +0.82%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.436"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("E36").Value = "  -0.68%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.106"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -2.76%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01964"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("E39").Value = "  -1.84%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.265"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.22%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.866"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  +0.93%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5028"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("E44").Value = "  -2.98%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.54"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -1.06%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.36"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -2.69%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4733"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  -0.66%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06389"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.86%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.661"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -2.17%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.880"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +0.03%  "
